# PopulationParameters.xlsx edit
#
# Per the commit: the "Demographics" sheet (previously the generic
# German default name "Tabelle1") now explicitly defines a small test
# population ("TestPopulation") used by runScenarios() for population
# simulations. The only real data change is the population size
# (numberOfIndividuals, column D) dropping from the placeholder 1000
# down to 2 for the lightweight test fixture, and the sheet is renamed
# to reflect its purpose. The previously-active "UserDefinedVariability"
# tab selection moves to "Demographics" (it's now the sheet users look
# at first), with the cell cursor left on D3 (just under the edited
# value).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: Tabelle1 -> Demographics ---------------------------------
$demographics = $wb.Worksheets.Item("Tabelle1")
$demographics.Name = "Demographics"

# numberOfIndividuals (column D, row 2) : 1000 -> 2
$demographics.Range("D2").Value = 2

# Make Demographics the active/selected sheet, cursor on D3.
$demographics.Activate() | Out-Null
$demographics.Range("D3").Select() | Out-Null
